# Generate Report for Handback
# This script updates the handback-status report so that the two source
# files tracked in the workbook now point at their newly generated GUIDs
# (and the merged/regenerated handoff-handback xlf pairs + timestamps).

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Old / new identifiers
# ----------------------------------------------------------------------
$oldGuid1 = "0d2bcf84-811f-4583-9f7c-c4eaf71dd853"
$newGuid1 = "bf25903d-bfde-49bb-a091-e574762c400c"
$oldGuid2 = "b1d04f01-21ec-4f04-813e-5af67ab12f26"
$newGuid2 = "ffff3a9f9841-4b06-4a6a-aea1-5755403f9426"

$newHash = "314048f49f06a58ff9c768703cf2246ce24a91ca"

$newMd1 = "$newGuid1.md"
$newMd2 = "$newGuid2.md"

$newXlfZh = "$newGuid1.$newHash.zh-cn.xlf"
$newXlfDe = "$newGuid1.$newHash.de-de.xlf"

$newHandoffZh1 = "2016-03-23 11:12:46"
$newHandbackZh1 = "2016-03-23 11:13:18"
$newHandoffDe1 = "2016-03-23 11:12:50"
$newHandbackDe1 = "2016-03-23 11:13:25"

# The diff only touches sharedStrings.xml and the sheetN.xml bodies; the
# worksheet .rels (i.e. the actual hyperlink target URLs) are left
# untouched, so we keep reusing the original (still GUID-1 / GUID-2
# based) URLs and only change the visible display text / cell value.

# ----------------------------------------------------------------------
# Sheet 1: Overview
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$addrA2 = "https://github.com/OpenLocalizationTest/oltest/blob/a7b86c6743aaa31f178806ee3b39c872899f4b34/e2e/$oldGuid1.md"
$addrA3 = "https://github.com/OpenLocalizationTest/oltest/blob/a7b86c6743aaa31f178806ee3b39c872899f4b34/e2e/$oldGuid2.md"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $addrA2, "", "", $newMd1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $addrA3, "", "", $newMd2)

# ----------------------------------------------------------------------
# Sheet 2: zh-cn
# ----------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhA2 = "https://github.com/OpenLocalizationTest/oltest/blob/a7b86c6743aaa31f178806ee3b39c872899f4b34/e2e/$oldGuid1.md"
$zhD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1972a376b49f255c11507b83853bd9afab7edeb9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid1.6e0d51801b200913843407a05ca15a8dcfb8d6ba.zh-cn.xlf"
$zhF2 = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b2c02e124ad76b57e8359cb8c013d4c3ad420c62/e2e/$oldGuid1.md"
$zhG2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/97a55f3749f1a7dca67689d8be0e389a1d90bb80/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid1.6e0d51801b200913843407a05ca15a8dcfb8d6ba.zh-cn.xlf"

$zhA3 = "https://github.com/OpenLocalizationTest/oltest/blob/a7b86c6743aaa31f178806ee3b39c872899f4b34/e2e/$oldGuid2.md"
$zhD3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1972a376b49f255c11507b83853bd9afab7edeb9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid2.c455b3698f1ca4287c6c1b886d0424549af4fb25.zh-cn.xlf"
$zhF3 = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/b2c02e124ad76b57e8359cb8c013d4c3ad420c62/e2e/$oldGuid2.md"
$zhG3 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/97a55f3749f1a7dca67689d8be0e389a1d90bb80/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid2.c455b3698f1ca4287c6c1b886d0424549af4fb25.zh-cn.xlf"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhA2, "", "", $newMd1)
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhD2, "", "", $newXlfZh)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhF2, "", "", $newMd1)
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhG2, "", "", $newXlfZh)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zhA3, "", "", $newMd2)
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhD3, "", "", $newXlfZh)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhF3, "", "", $newMd2)
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhG3, "", "", $newXlfZh)

$wsZh.Range("E2").Value = $newHandoffZh1
$wsZh.Range("H2").Value = $newHandbackZh1
$wsZh.Range("E3").Value = $newHandoffZh1
$wsZh.Range("H3").Value = $newHandbackZh1

# ----------------------------------------------------------------------
# Sheet 3: de-de
# ----------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deA2 = "https://github.com/OpenLocalizationTest/oltest/blob/a7b86c6743aaa31f178806ee3b39c872899f4b34/e2e/$oldGuid1.md"
$deD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a7bde0d5ead474229e22b9a05e381482d8511c09/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid1.6e0d51801b200913843407a05ca15a8dcfb8d6ba.de-de.xlf"
$deF2 = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/962d8883636b1d28f44e089a4d4a5695ab5abfd2/e2e/$oldGuid1.md"
$deG2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/95f905ff7ea8e88288ab0d48595ec82f76b84ac0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid1.6e0d51801b200913843407a05ca15a8dcfb8d6ba.de-de.xlf"

$deA3 = "https://github.com/OpenLocalizationTest/oltest/blob/a7b86c6743aaa31f178806ee3b39c872899f4b34/e2e/$oldGuid2.md"
$deD3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a7bde0d5ead474229e22b9a05e381482d8511c09/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid2.c455b3698f1ca4287c6c1b886d0424549af4fb25.de-de.xlf"
$deF3 = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/962d8883636b1d28f44e089a4d4a5695ab5abfd2/e2e/$oldGuid2.md"
$deG3 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/95f905ff7ea8e88288ab0d48595ec82f76b84ac0/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid2.c455b3698f1ca4287c6c1b886d0424549af4fb25.de-de.xlf"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deA2, "", "", $newMd1)
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deD2, "", "", $newXlfDe)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deF2, "", "", $newMd1)
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deG2, "", "", $newXlfDe)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $deA3, "", "", $newMd2)
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deD3, "", "", $newXlfDe)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deF3, "", "", $newMd2)
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deG3, "", "", $newXlfDe)

$wsDe.Range("E2").Value = $newHandoffDe1
$wsDe.Range("H2").Value = $newHandbackDe1
$wsDe.Range("E3").Value = $newHandoffDe1
$wsDe.Range("H3").Value = $newHandbackDe1
